$d = $word.ActiveDocument

$replacements = @(
    @("667÷6=", "878÷8="),
    @("909÷2=", "827÷2="),
    @("226÷4=", "596÷4="),
    @("241÷8=", "856÷5="),
    @("231÷4=", "105÷5="),
    @("216÷4=", "365÷9="),
    @("992÷6=", "541÷6="),
    @("467÷5=", "924÷8="),
    @("750÷5=", "464÷5="),
    @("663÷2=", "552÷8="),
    @("170÷8=", "368÷5="),
    @("419÷9=", "733÷4="),
    @("552÷5=", "566÷3="),
    @("130÷2=", "371÷9="),
    @("863÷5=", "803÷2="),
    @("301÷8=", "576÷2="),
    @("206÷9=", "407÷5="),
    @("564÷2=", "295÷2="),
    @("711÷2=", "791÷4="),
    @("483÷4=", "712÷2="),
    @("354÷6=", "784÷4="),
    @("467÷8=", "468÷8="),
    @("297÷4=", "795÷3="),
    @("652÷3=", "665÷3="),
    @("876÷2=", "606÷7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
